$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: "34172" (number, text-formatted) -> "34712" (plain text, General format)
$ws.Range("C2").Value = "34712"
$ws.Range("C2").NumberFormat = "General"

# J2 / L2 / M2: boolean constants -> TRUE()/FALSE() formulas, displayed as text ("@")
$ws.Range("J2").Formula = "=TRUE()"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("L2").Formula = "=FALSE()"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("M2").Formula = "=TRUE()"
$ws.Range("M2").NumberFormat = "@"

# Move the active selection to C2 (was M31)
$ws.Range("C2").Select()
